$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: fill in new LeetCode #13 "Roman To Integer" entry (Anna) ---
$ws.Range("D37").WrapText = $true
$ws.Range("D37").Value = "13. Roman To Integer"
$ws.Range("A37").EntireRow.RowHeight = 14

# --- Row 38: same entry for Stephan, including date + category ---
$ws.Range("D38").WrapText = $true
$ws.Range("D38").Value = "13. Roman To Integer"
$ws.Range("E38").Value = "2020/12/28"
$ws.Range("G38").Value = "Completed"
$ws.Range("A38").EntireRow.RowHeight = 14

# --- Row 39: new blank template row (Anna) with styled-but-empty D cell ---
$ws.Range("A39").Value = "LeetCode"
$ws.Range("B39").Value = "Anna"
$ws.Range("C39").Value = "Easy"
$ws.Range("D39").WrapText = $true

# --- Row 40: new blank template row (Stephan) ---
$ws.Range("A40").Value = "LeetCode"
$ws.Range("B40").Value = "Stephan"
$ws.Range("C40").Value = "Easy"

# --- Row 41: new blank template row (Anna) ---
$ws.Range("A41").Value = "LeetCode"
$ws.Range("B41").Value = "Anna"
$ws.Range("C41").Value = "Easy"

# --- Row 42: new blank template row (Stephan) ---
$ws.Range("A42").Value = "LeetCode"
$ws.Range("B42").Value = "Stephan"
$ws.Range("C42").Value = "Easy"

# Match the alignment style used by the other rows' A/B columns ("Site"/"User")
$ws.Range("A39:A42").HorizontalAlignment = -4131
$ws.Range("B39:B42").HorizontalAlignment = -4131

# --- Move selection to D39 (matches the saved cursor position in the edit) ---
$ws.Range("D39").Select() | Out-Null
